$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 18 is currently an empty/unused row (a gap in the data).
# Populate it with a new entry "L_LongI", which mirrors the row above it (L_I)
# but represents a new dimension parameter.
$ws.Range("A18").Value = "L_LongI"
$ws.Range("B18").Value = 35.200000000000003
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 2.4500000000000002
$ws.Range("E18").Value = 17.600000000000001
$ws.Range("F18").Value = 0.20300000000000001

# Update the view: scroll position and the active selection.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("A18").Select()
